$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dSF column (F) values per repulled data / mean calculation
$ws.Range("F2").Value = -4
$ws.Range("F5").Value = -5
$ws.Range("F7").Value = -2
$ws.Range("F8").Value = -4
$ws.Range("F13").Value = -13
$ws.Range("F16").Value = -3
$ws.Range("F19").Value = -5
